$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.352.07'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '1.787.28'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '''226.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").Value = '''0.555'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '''32.87'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.65%  '
$ws.Range("E9").Value = '  +1.31%  '
$ws.Range("D10").Value = '''0.0689'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.72%  '
$ws.Range("D11").Value = '''0.0946'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '2.044.63'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("D13").Value = '''11.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.88%  '
$ws.Range("D14").Value = '1.795.73'
$ws.Range("E14").Value = '  +0.76%  '
$ws.Range("E15").Value = '  +2.31%  '
$ws.Range("D16").Value = '34.336.66'
$ws.Range("E16").Value = '  +0.76%  '
$ws.Range("E17").Value = '  +2.72%  '
$ws.Range("D18").Value = '''68.38'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("D19").Value = '''244.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("D20").Value = '0.0₃0794'
$ws.Range("E20").Value = '  +1.07%  '
$ws.Range("D21").Value = '''11.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.73%  '
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("D24").Value = '''168.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.48%  '
$ws.Range("E25").Value = '  +2.06%  '
$ws.Range("D26").Value = '''7.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.38%  '
$ws.Range("D27").Value = '''16.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.86%  '
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").Value = '''4.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.92%  '
$ws.Range("D31").Value = '''0.0527'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.15%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''1.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''3.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("E34").Value = '  +1.28%  '
$ws.Range("D35").Value = '1.411.71'
$ws.Range("E35").Value = '  -2.14%  '
$ws.Range("D36").Value = '''2.58'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.57%  '
$ws.Range("E37").Value = '  +5.04%  '
$ws.Range("E38").Value = '  +2.95%  '
$ws.Range("E39").Value = '  +0.49%  '
$ws.Range("D40").Value = '''84.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.83%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''14.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.00%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").Value = '''2.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").Value = '''2.78'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.49%  '
$ws.Range("D44").Value = '''0.938'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.80%  '
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("E46").Value = '  +2.70%  '
$ws.Range("D47").Value = '''6.07'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("D48").Value = '1.946.15'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").Value = '''105.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("E51").Value = '  -1.94%  '
